# Weekly update: prepend a new week's worth of data (2 rows) at row 145,
# pushing the existing rows 145:172 down to 147:174.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 145 (shifts 145:172 -> 147:174)
$ws.Rows.Item(145).Resize(2).Insert()

# --- New row 145 ---
$ws.Range("A145").Value = 9
$ws.Range("B145").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C145").Value = "Metropolitana"
$ws.Range("D145").Value = 44511
$ws.Range("E145").Value = 13
$ws.Range("F145").Value = 100112043
$ws.Range("G145").Value = "Pepino ensalada"
$ws.Range("H145").Value = "Sin especificar"
$ws.Range("I145").Value = "Primera"
$ws.Range("J145").Value = 106
$ws.Range("K145").Value = 7000
$ws.Range("L145").Value = 8000
$ws.Range("M145").Value = 7500
$ws.Range("N145").Value = "`$/caja 60 unidades"
$ws.Range("O145").Value = "Región de Arica y Parinacota"
$ws.Range("P145").Value = 125
$ws.Range("Q145").Value = 60
$ws.Range("R145").Value = "Hortaliza"

# --- New row 146 ---
$ws.Range("A146").Value = 9
$ws.Range("B146").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C146").Value = "Metropolitana"
$ws.Range("D146").Value = 44511
$ws.Range("E146").Value = 13
$ws.Range("F146").Value = 100112043
$ws.Range("G146").Value = "Pepino ensalada"
$ws.Range("H146").Value = "Sin especificar"
$ws.Range("I146").Value = "Segunda"
$ws.Range("J146").Value = 61
$ws.Range("K146").Value = 6000
$ws.Range("L146").Value = 6000
$ws.Range("M146").Value = 6000
$ws.Range("N146").Value = "`$/caja 100 unidades"
$ws.Range("O146").Value = "Región de Arica y Parinacota"
$ws.Range("P146").Value = 60
$ws.Range("Q146").Value = 100
$ws.Range("R146").Value = "Hortaliza"
